$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell contents (headers + data rows) ---
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Result"

$ws.Range("A2").Value = "catdog"
$ws.Range("B2").Value = "catdog"
$ws.Range("C2").Value = "Pass"

$ws.Range("A3").Value = "dogcat"
$ws.Range("B3").Value = "dogcat"
$ws.Range("C3").Value = "Pass"

$ws.Range("A4").Value = "dogcat"
$ws.Range("B4").Value = "catdog"
$ws.Range("C4").Value = "Pass"

$ws.Range("A5").Value = "catdog"
$ws.Range("B5").Value = "dogcat"
$ws.Range("C5").Value = "Pass"

# --- Column widths: drop custom widths on columns B and C, set column A to 17 ---
$ws.Columns.Item(2).EntireColumn.ClearFormats()
$ws.Columns.Item(3).EntireColumn.ClearFormats()
$ws.Columns.Item(1).ColumnWidth = 16.15

# --- Sheet view: zoom to 190% and move the selection to E8 ---
$excel.ActiveWindow.Zoom = 190
[void]$ws.Range("E8").Select()
